$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (price/volume refresh + Stellar/Cosmos row swap).
# Each value is written with a leading apostrophe so that numeric-looking
# text (e.g. "27.232.18", "7.31") is stored verbatim as text rather than
# being reinterpreted by Excel as a number; the Style is then reset to
# "Normal" so no stray quote-prefix / text-format style is left behind.
$ws.Range('D2').Value = "'" + '27.232.18'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + '  +0.04%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + '1.635.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + '  -0.96%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + '  -0.06%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + '216.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + '  -0.41%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + '0.522'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + '  +0.84%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + '  -0.03%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'" + '0.256'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'" + '  -0.45%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'" + '  -0.36%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'" + '20.43'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'" + '  +1.87%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'" + '0.0850'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'" + '  -0.10%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + '1.648.68'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + '  -0.31%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + '4.14'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + '  -0.28%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + '0.547'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + '  +0.81%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + '65.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + '  -3.46%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'" + '27.196.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'" + '  -0.11%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').Value = "'" + '  +0.20%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + '218.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + '  -0.64%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'" + '  -0.04%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'" + '6.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + '  +1.82%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + '4.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + '  -0.83%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'" + '2.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'" + '  -5.77%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + '9.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + '  -1.51%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'" + '147.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'" + '  +0.05%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'" + '  -0.11%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = "'" + 'Cosmos'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'" + 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'" + '7.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + '  -3.22%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('B27').Value = "'" + 'Stellar'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = "'" + '0.118'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + '  -0.03%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'" + '15.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'" + '  -1.14%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'" + '  -0.59%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'" + '  -0.58%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'" + '  -0.42%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + '  -1.27%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + '1.322.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + '  +4.02%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'" + '  -0.35%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'" + '  -0.47%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'" + '0.0176'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'" + '  -0.96%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + '0.547'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + '  -0.24%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'" + '0.851'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'" + '  +0.30%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'" + '  -0.07%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'" + '  +1.60%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + '0.803'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + '  -0.82%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'" + '64.44'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'" + '  +3.66%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'" + '1.773.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'" + '  -1.05%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'" + '5.24'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'" + '  -3.87%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + '91.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + '  -0.75%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + '1.62'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + '  +1.04%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + '0.0₆0105'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + '  +0.15%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + '0.808'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + '  +21.28%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + '  +0.18%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'" + '0.0989'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'" + '  +1.53%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'" + '7.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'" + '  -1.21%  '
$ws.Range('E51').Style = 'Normal'
